$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new columns (H, I) for TexturePacker / UIAtlas support
$ws.Range("H1").Value = "simple_atlas"
$ws.Range("I1").Value = "simple_tex"

$ws.Range("H2").Value = "string"
$ws.Range("I2").Value = "string"

$ws.Range("H3").Value = "缩略图集"
$ws.Range("I3").Value = "缩略图"

$ws.Range("H4").Value = "CardSimple"

$ws.Range("J11").Select()
